$d = $word.ActiveDocument

# --- wdReplaceAll / wdFindContinue constants used below ---
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
#   Wrap: 1 = wdFindContinue, Replace: 2 = wdReplaceAll

# 1) Title
$d.Content.Find.Execute("Unraveling the Enigma of Neuronal Connectivity", $true, $false, $false, $false, $false, $true, 1, $false, "Chemistry: A World of Molecules and Reactions", 2)

# 2) Author name (keep "Dr." prefix run untouched, only replace the name run's text)
$d.Content.Find.Execute(" Naomi Williams", $true, $false, $false, $false, $false, $true, 1, $false, " Kimberly Young", 2)

# 3) Email line: collapse "williams" + "." + "naomi@neuroscience" runs into a single new run,
#    leaving the trailing "." and "edu" runs untouched.
$d.Content.Find.Execute("williams.naomi@neuroscience", $true, $false, $false, $false, $false, $true, 1, $false, "kyyoung@highlandschools", 2)

# 4) Body paragraph sentences (each sentence/run replaced in place; separating "." runs are left as-is)
$d.Content.Find.Execute("Within the intricate tapestry of the human brain, a universe of neurons forms a dynamic network, orchestrating our every thought, action, and emotion", $true, $false, $false, $false, $false, $true, 1, $false, "Chemistry, an integral field of science, explores the nature of matter and its interactions", 2)

$d.Content.Find.Execute(" Understanding the intricate web of connections between these neurons, known as neuronal connectivity, holds the key to unraveling the mysteries of the mind", $true, $false, $false, $false, $false, $true, 1, $false, " Delving into this realm unveils the fundamental building blocks of the universe: molecules and atoms", 2)

$d.Content.Find.Execute(" As we embark on this journey of exploration, we delve into the world of neuroscience, where scientists diligently decipher the language of the brain, one synapse at a time", $true, $false, $false, $false, $false, $true, 1, $false, " Chemistry provides a window into the intricate world of chemical reactions, where substances undergo transformations, resulting in the formation of new substances with distinct properties. Through its principles and applications, chemistry plays a pivotal role in understanding various phenomena observed in the natural world", 2)

$d.Content.Find.Execute("From the bustling metropolis of New York City to the serene countryside of Provence, researchers toil tirelessly in their laboratories, meticulously tracing the intricate pathways of neuronal communication", $true, $false, $false, $false, $false, $true, 1, $false, "Unveiling the fundamental principles of chemistry enables us to unravel the secrets of molecular behavior", 2)

$d.Content.Find.Execute(" They employ a symphony of cutting-edge technologies, from high-resolution microscopes that peer into the depths of the brain to computational algorithms that analyze vast troves of data", $true, $false, $false, $false, $false, $true, 1, $false, " The periodic table, a cornerstone of chemistry, organizes elements based on their atomic number, unveiling periodic trends that govern their properties and reactivities", 2)

$d.Content.Find.Execute(" Each breakthrough, each discovery, brings us closer to comprehending the enigmatic dance of neurons that underpins our existence", $true, $false, $false, $false, $false, $true, 1, $false, " By mastering these principles, we unlock the ability to predict and manipulate chemical reactions, paving the way for countless innovations and technological advancements that shape our modern world", 2)

$d.Content.Find.Execute("The quest to understand neuronal connectivity is not merely an academic pursuit; it holds immense promise for unraveling neurological and psychiatric disorders that afflict millions worldwide", $true, $false, $false, $false, $false, $true, 1, $false, "Chemistry's far-reaching impact extends beyond the laboratory walls", 2)

$d.Content.Find.Execute(" By deciphering the intricate patterns of neuronal communication, we can illuminate the disruptions that give rise to conditions such as Alzheimer's disease, schizophrenia, and autism", $true, $false, $false, $false, $false, $true, 1, $false, " It underpins the development of medicines that alleviate human suffering, fuels the engines that power our transportation, and enables the creation of materials that enhance our daily lives", 2)

$d.Content.Find.Execute(" Armed with this knowledge, we can pave the way for targeted therapies that restore the delicate balance of neuronal connectivity, alleviating suffering and restoring hope", $true, $false, $false, $false, $false, $true, 1, $false, " Its applications encompass agriculture, energy production, and environmental protection, showcasing the diverse contributions chemistry makes to society", 2)

# 5) Summary heading paragraph body
$d.Content.Find.Execute("In this essay, we have embarked on a journey into the realm of neuronal connectivity, exploring the intricate network of connections that orchestrates the symphony of our thoughts, actions, and emotions", $true, $false, $false, $false, $false, $true, 1, $false, "Venturing into the captivating realm of chemistry, we discover the fundamental principles that govern the nature of matter and its interactions", 2)

$d.Content.Find.Execute(" We have delved into the world of neuroscience, where scientists tirelessly decipher the language of the brain, employing cutting-edge technologies and computational algorithms to unravel the mysteries of neuronal communication", $true, $false, $false, $false, $false, $true, 1, $false, " By unveiling the secrets of molecular behavior and chemical reactions, chemistry enables us to understand various phenomena in the world around us", 2)

# This also removes the embedded <w:lastRenderedPageBreak/> that used to sit between the two runs being merged here.
$d.Content.Find.Execute(" Our understanding of neuronal connectivity holds immense promise for illuminating neurological and psychiatric disorders, paving the way for targeted therapies that restore the delicate balance of neuronal connectivity and alleviate suffering", $true, $false, $false, $false, $false, $true, 1, $false, " Its applications permeate numerous fields, from medicine and energy production to material science and environmental protection, underscoring its vital role in shaping our world. Chemistry continues to inspire and challenge, offering boundless possibilities for exploration and innovation", 2)

# 6) Add a new empty paragraph at the very end of the document body.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
